$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B3 date (was showing 2026-06-26, corrected to 2024-06-26) -- keep existing style.
$ws.Range("B3").Value = 45469

# Row 4
$ws.Range("A4").Formula = "'03"
$ws.Range("C4").Value = "Selección ENDI año 3"
$ws.Range("D4").Value = "Angel Gaibor"

# Row 5
$ws.Range("A5").Value = "04"
$ws.Range("C5").Value = "Selección ENCIET 202410"
$ws.Range("D5").Value = "Rafael Encalada"

# Row 6 (accion entered before version, matching original author's order)
$ws.Range("C6").Value = "Selección ENCIET 202411"
$ws.Range("A6").Value = "05"
$ws.Range("D6").Value = "Rafael Encalada"

# Row 7
$ws.Range("A7").Value = "06"
$ws.Range("C7").Value = "Selección ENCIET 202412"
$ws.Range("D7").Value = "Rafael Encalada"

# Row 8
$ws.Range("A8").Formula = "'07"
$ws.Range("C8").Value = "Selección ENDI año 3 6 viv"
$ws.Range("D8").Value = "Angel Gaibor"

# Row 9
$ws.Range("A9").Formula = "'08"
$ws.Range("C9").Value = "Selección ENDI año 3 7 viv"
$ws.Range("D9").Value = "Angel Gaibor"

# Row 10
$ws.Range("A10").Value = "09"
$ws.Range("B10").Value = 45671
$ws.Range("B10").NumberFormat = "mm-dd-yy"
$ws.Range("C10").Value = "Borrado selección ENDI 3 muestras"
$ws.Range("D10").Value = "Angel Gaibor"

$ws.Range("B4").Select() | Out-Null
